# Change Singleton presentation footer color
#
# For every slide that has Date / Footer / Slide Number placeholders:
#   - colour the Date placeholder's text (the datetime field run) green
#   - colour the Slide Number placeholder's text green
#   - remove the "@Bellkross" Footer placeholder entirely
#
# ppPlaceholderDate = 16, ppPlaceholderFooter = 15, ppPlaceholderSlideNumber = 13
# Target colour: RGB(167, 216, 109) == hex A7D86D == 167 + 216*256 + 109*65536

$greenRgb = 167 + (216 * 256) + (109 * 65536)

$p = $ppt.ActivePresentation

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)

    # Walk shapes back-to-front since we may delete some along the way.
    for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
        $shp = $s.Shapes.Item($i)

        if (-not $shp.HasTextFrame) {
            continue
        }
        if ($shp.PlaceholderFormat -eq $null) {
            continue
        }

        $phType = 0
        try {
            $phType = $shp.PlaceholderFormat.Type
        } catch {
            continue
        }

        if ($phType -eq 16) {
            # Date placeholder -> green text
            $shp.TextFrame.TextRange.Font.Color.RGB = $greenRgb
        }
        elseif ($phType -eq 13) {
            # Slide Number placeholder -> green text
            $shp.TextFrame.TextRange.Font.Color.RGB = $greenRgb
        }
    }

    # Turn the footer off for this slide; PowerPoint removes the
    # "@Bellkross" footer placeholder shape from the slide XML.
    $hf = $s.HeadersFooters
    if ($hf.Footer.Visible) {
        $hf.Footer.Visible = $false
    }
}
